# Applies the cryptos.xlsx price/volume refresh described in the commit.
# Cells are plain text in the source workbook (t="inlineStr"); Excel's COM
# layer auto-detects numeric-looking strings and would silently convert them
# to numbers, so every write is bracketed with a temporary "@" (Text) number
# format and then restored to the default "Normal" style, matching the
# original (unstyled) cells exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "60.834.75"
Set-TextValue "E2" "  +0.00%  "
Set-TextValue "D3" "2.378.60"
Set-TextValue "E3" "  -3.12%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "544.52"
Set-TextValue "E5" "  -0.76%  "
Set-TextValue "D6" "141.44"
Set-TextValue "E6" "  -2.56%  "
Set-TextValue "E7" "  +0.10%  "
Set-TextValue "D8" "0.547"
Set-TextValue "E8" "  -8.24%  "
Set-TextValue "D9" "2.377.60"
Set-TextValue "E9" "  -3.08%  "
Set-TextValue "E10" "  -1.47%  "
Set-TextValue "E11" "  +0.79%  "
Set-TextValue "E12" "  -0.62%  "
Set-TextValue "E13" "  -1.63%  "
Set-TextValue "D14" "25.54"
Set-TextValue "E14" "  -1.74%  "
Set-TextValue "D15" "2.809.22"
Set-TextValue "E15" "  -2.91%  "
Set-TextValue "E16" "  -0.34%  "
Set-TextValue "D17" "60.625.54"
Set-TextValue "E17" "  -0.19%  "
Set-TextValue "D18" "2.379.94"
Set-TextValue "E18" "  -3.01%  "
Set-TextValue "D19" "10.65"
Set-TextValue "E19" "  -3.53%  "
Set-TextValue "D20" "4.10"
Set-TextValue "E20" "  -1.53%  "
Set-TextValue "D21" "316.60"
Set-TextValue "D22" "6.71"
Set-TextValue "E22" "  -2.73%  "
Set-TextValue "E23" "  -0.07%  "
Set-TextValue "D24" "1.84"
Set-TextValue "E24" "  +4.32%  "
Set-TextValue "D25" "62.86"
Set-TextValue "E25" "  -0.81%  "
Set-TextValue "E26" "  +0.06%  "
Set-TextValue "D27" "2.496.19"
Set-TextValue "E27" "  -3.06%  "
Set-TextValue "D28" "0.0₃0932"
Set-TextValue "E28" "  -4.91%  "
Set-TextValue "D29" "7.77"
Set-TextValue "E29" "  +2.30%  "
Set-TextValue "D30" "521.40"
Set-TextValue "E30" "  -2.93%  "
Set-TextValue "E31" "  -4.13%  "
Set-TextValue "D32" "8.00"
Set-TextValue "E32" "  -3.90%  "
Set-TextValue "E33" "  -3.86%  "
Set-TextValue "E34" "  -2.85%  "
Set-TextValue "E35" "  -0.50%  "
Set-TextValue "E36" "  +0.04%  "
Set-TextValue "E37" "  -6.28%  "
Set-TextValue "D38" "4.67"
Set-TextValue "E39" "  -0.40%  "
Set-TextValue "E40" "  -2.16%  "
Set-TextValue "E41" "  +1.11%  "
Set-TextValue "E42" "  +0.12%  "
Set-TextValue "D43" "137.43"
Set-TextValue "E43" "  -4.97%  "
Set-TextValue "D44" "40.31"
Set-TextValue "E44" "  +1.41%  "
Set-TextValue "E45" "  -2.77%  "
Set-TextValue "D46" "139.91"
Set-TextValue "E46" "  -4.48%  "
Set-TextValue "D47" "3.56"
Set-TextValue "E47" "  +0.13%  "
Set-TextValue "D48" "20.36"
Set-TextValue "E48" "  -2.18%  "
Set-TextValue "D49" "0.0519"
Set-TextValue "E49" "  -1.88%  "
Set-TextValue "E50" "  -1.31%  "
Set-TextValue "E51" "  -2.67%  "

Write-Host "Applied 76 cell updates"
